$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.737.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.738.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.736.06"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.79%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E10").Value = "  +2.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.361.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.738.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.726.05"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.30"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +15.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  -2.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -3.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.98"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.882.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.670.71"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.80%  "

$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.93"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.19"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.806.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("E51").Value = "  +0.68%  "

